$d = $word.ActiveDocument

$d.Content.Find.Execute("68×93=6324", $true, $false, $false, $false, $false, $true, 1, $false, "60×66=3960", 2)
$d.Content.Find.Execute("22×97=2134", $true, $false, $false, $false, $false, $true, 1, $false, "73×81=5913", 2)
$d.Content.Find.Execute("51×61=3111", $true, $false, $false, $false, $false, $true, 1, $false, "99×83=8217", 2)
$d.Content.Find.Execute("19×44=836", $true, $false, $false, $false, $false, $true, 1, $false, "19×28=532", 2)
$d.Content.Find.Execute("73×49=3577", $true, $false, $false, $false, $false, $true, 1, $false, "61×53=3233", 2)
$d.Content.Find.Execute("62×45=2790", $true, $false, $false, $false, $false, $true, 1, $false, "72×38=2736", 2)
$d.Content.Find.Execute("96×99=9504", $true, $false, $false, $false, $false, $true, 1, $false, "86×26=2236", 2)
$d.Content.Find.Execute("60×35=2100", $true, $false, $false, $false, $false, $true, 1, $false, "40×22=880", 2)
$d.Content.Find.Execute("31×61=1891", $true, $false, $false, $false, $false, $true, 1, $false, "48×16=768", 2)
$d.Content.Find.Execute("34×77=2618", $true, $false, $false, $false, $false, $true, 1, $false, "79×73=5767", 2)
$d.Content.Find.Execute("24×57=1368", $true, $false, $false, $false, $false, $true, 1, $false, "98×38=3724", 2)
$d.Content.Find.Execute("46×35=1610", $true, $false, $false, $false, $false, $true, 1, $false, "81×57=4617", 2)
$d.Content.Find.Execute("50×63=3150", $true, $false, $false, $false, $false, $true, 1, $false, "43×99=4257", 2)
$d.Content.Find.Execute("41×92=3772", $true, $false, $false, $false, $false, $true, 1, $false, "42×28=1176", 2)
$d.Content.Find.Execute("25×56=1400", $true, $false, $false, $false, $false, $true, 1, $false, "44×82=3608", 2)
$d.Content.Find.Execute("73×29=2117", $true, $false, $false, $false, $false, $true, 1, $false, "71×11=781", 2)
$d.Content.Find.Execute("31×90=2790", $true, $false, $false, $false, $false, $true, 1, $false, "79×33=2607", 2)
$d.Content.Find.Execute("90×27=2430", $true, $false, $false, $false, $false, $true, 1, $false, "62×39=2418", 2)
$d.Content.Find.Execute("25×13=325", $true, $false, $false, $false, $false, $true, 1, $false, "73×80=5840", 2)
$d.Content.Find.Execute("18×88=1584", $true, $false, $false, $false, $false, $true, 1, $false, "39×90=3510", 2)
$d.Content.Find.Execute("85×91=7735", $true, $false, $false, $false, $false, $true, 1, $false, "51×50=2550", 2)
$d.Content.Find.Execute("25×87=2175", $true, $false, $false, $false, $false, $true, 1, $false, "89×79=7031", 2)
$d.Content.Find.Execute("63×86=5418", $true, $false, $false, $false, $false, $true, 1, $false, "39×91=3549", 2)
$d.Content.Find.Execute("69×57=3933", $true, $false, $false, $false, $false, $true, 1, $false, "41×96=3936", 2)
$d.Content.Find.Execute("61×83=5063", $true, $false, $false, $false, $false, $true, 1, $false, "13×49=637", 2)
